$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values
$ws.Range("B3").Value = 5.5
$ws.Range("C4").Value = 1.25

# Update selection to match the saved UI state
$ws.Range("C4").Select()
